$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold+border) from Z1 into the new header cells AA1:AD1
$ws.Range("Z1").Copy()
$ws.Range("AA1:AD1").PasteSpecial(-4122)

# Clear existing contents first so that the shared-string table gets rebuilt
# strictly in the row-major cell order used below (matches target layout).
$ws.Range("B1:Z3").ClearContents()

# Row 1 headers
$ws.Range("B1").Value = "Unnamed: 0"
$ws.Range("C1").Value = "MI filename"
$ws.Range("D1").Value = "Water_to_HW_ratio_Trapezoid"
$ws.Range("E1").Value = "Water_to_HW_ratio_Simpson"
$ws.Range("F1").Value = "Water_to_Total_Silicate_ratio_Trapezoid"
$ws.Range("G1").Value = "Water_to_Total_Silicate_ratio_Simpson"
$ws.Range("H1").Value = "Water_Trapezoid_Area"
$ws.Range("I1").Value = "Water_Simpson_Area"
$ws.Range("J1").Value = "Silicate_Trapezoid_Area"
$ws.Range("K1").Value = "Silicate_Simpson_Area"
$ws.Range("L1").Value = "Silicate_LHS_Back1"
$ws.Range("M1").Value = "Silicate_LHS_Back2"
$ws.Range("N1").Value = "Silicate_RHS_Back1"
$ws.Range("O1").Value = "Silicate_RHS_Back2"
$ws.Range("P1").Value = "Silicate_N_Poly"
$ws.Range("Q1").Value = "LW_Silicate_Trapezoid_Area"
$ws.Range("R1").Value = "LW_Silicate_Simpson_Area"
$ws.Range("S1").Value = "HW_Silicate_Trapezoid_Area"
$ws.Range("T1").Value = "HW_Silicate_Simpson_Area"
$ws.Range("U1").Value = "MW_Silicate_Trapezoid_Area"
$ws.Range("V1").Value = "MW_Silicate_Simpson_Area"
$ws.Range("W1").Value = "Water Filename"
$ws.Range("X1").Value = "Water_LHS_Back1"
$ws.Range("Y1").Value = "Water_LHS_Back2"
$ws.Range("Z1").Value = "Water_RHS_Back1"
$ws.Range("AA1").Value = "Water_RHS_Back2"
$ws.Range("AB1").Value = "Water_N_Poly"
$ws.Range("AC1").Value = "HW:LW_Trapezoid"
$ws.Range("AD1").Value = "HW:LW_Simpson"

# Row 2 data
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = "ETFS_OL39_MI7_50X_GLASS.txt"
$ws.Range("D2").Value = 1.035972117635653
$ws.Range("E2").Value = 1.036197547193371
$ws.Range("F2").Value = 0.7579678818556516
$ws.Range("G2").Value = 0.7576781215089602
$ws.Range("H2").Value = 274807.3951902158
$ws.Range("I2").Value = 274611.0037276433
$ws.Range("J2").Value = 362558.1001103033
$ws.Range("K2").Value = 362437.5522164207
$ws.Range("L2").Value = 300
$ws.Range("M2").Value = 340
$ws.Range("N2").Value = 1200
$ws.Range("O2").Value = 1250
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 69943.66400954682
$ws.Range("R2").Value = 69904.77062547633
$ws.Range("S2").Value = 265265.2426760239
$ws.Range("T2").Value = 265018.0020899012
$ws.Range("U2").Value = 11333.3329334461
$ws.Range("V2").Value = 11224.56189991856
$ws.Range("W2").Value = "ETFS_OL39_MI7_50X_GLASS.txt"
$ws.Range("X2").Value = 2500
$ws.Range("Y2").Value = 2750
$ws.Range("Z2").Value = 3750
$ws.Range("AA2").Value = 4100
$ws.Range("AB2").Value = 3
$ws.Range("A2").Value = 0

# Row 3 data
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = "test_H2O.txt"
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("G3").ClearContents()
$ws.Range("H3").Value = 274807.3951902158
$ws.Range("I3").Value = 274611.0037276433
$ws.Range("J3").Value = 362558.1001103033
$ws.Range("K3").Value = 362437.5522164207
$ws.Range("L3").Value = 300
$ws.Range("M3").Value = 340
$ws.Range("N3").Value = 1200
$ws.Range("O3").Value = 1250
$ws.Range("P3").Value = 4
$ws.Range("Q3").Value = 69943.66400954682
$ws.Range("R3").Value = 69904.77062547633
$ws.Range("S3").Value = 69943.66400954682
$ws.Range("T3").Value = 69904.77062547633
$ws.Range("U3").Value = 11333.3329334461
$ws.Range("V3").Value = 11224.56189991856
$ws.Range("W3").Value = "test_H2O.txt"
$ws.Range("X3").Value = 2500
$ws.Range("Y3").Value = 2750
$ws.Range("Z3").Value = 3750
$ws.Range("AA3").Value = 4100
$ws.Range("AB3").Value = 3
$ws.Range("AC3").Value = 3.928981975446791
$ws.Range("AD3").Value = 3.928358555082121
$ws.Range("A3").Value = 1
